# Applies the "Added test methods in Profile test" commit to the
# TestScript-admin.xlsx workbook:
#   1. Flip the "Execute" column (B) from Yes -> No for every row that
#      was still set to "Yes" (the previously-enabled test rows).
#   2. Append three new test rows (64-66) for the new Profile/Api
#      Business test cases, with the final appended row left as the
#      one "Yes" (enabled) row, mirroring the others being turned off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Execute column: Yes -> No for the rows that were enabled ---
$rowsToDisable = @(12, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 32, 34, 35, 36)
foreach ($r in $rowsToDisable) {
    $ws.Cells.Item($r, 2).Value = "No"
}

# --- 2. Append the new rows describing the Api Business tests ---

# Row 64: create Api Business
$ws.Range("A64").Value = "create Api Business "
$ws.Range("B64").Value = "No"
$ws.Range("C64").Value = "testdata-admin.xlsx,profileApiBusiness"
$ws.Range("D64").Value = "RunOneIteration"
$ws.Range("E64").Value = "1"
$ws.Range("F64").Value = "1"
$ws.Range("G64").Value = "Api Business"
$ws.Range("H64").Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode"
$ws.Range("I64").Value = "coyni.admin.tests.ProfilesTest,`ntestAddApiBusiness,`n-pheadingInvitation,`n-pheadingContact,`n-pfirstName,`n-plastName,`n-pemail1,`n-pphoneNumber,`n-pbusiness,`n-pdoller,`n-ppercentage,`n-pecoSystem,`n-pbusinessLimit,`n-pecoSystemLimit"

# Row 65: Verify Business userList
$ws.Range("A65").Value = "Verify Business userList"
$ws.Range("B65").Value = "No"
$ws.Range("C65").Value = "testdata-admin.xlsx,profileApiBusiness"
$ws.Range("D65").Value = "RunOneIteration"
$ws.Range("E65").Value = "1"
$ws.Range("F65").Value = "1"
$ws.Range("G65").Value = "Api Business"
$ws.Range("H65").Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode"
$ws.Range("I65").Value = "coyni.admin.tests.ProfilesTest,`ntestveifyGrid,`n-pfilterType,`n-ptoAmount,`n-pamount,`n-pheading,`n-pexportHeading"
$ws.Range("I65").WrapText = $true
$ws.Range("I65").VerticalAlignment = -4160

# Row 66: Verify Business Details (left enabled = "Yes")
$ws.Range("A66").Value = "Verify Business Details"
$ws.Range("B66").Value = "Yes"
$ws.Range("C66").Value = "testdata-admin.xlsx,profileApiBusiness"
$ws.Range("D66").Value = "RunOneIteration"
$ws.Range("E66").Value = "1"
$ws.Range("F66").Value = "1"
$ws.Range("G66").Value = "Api Business"
$ws.Range("H66").Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode"
$ws.Range("I66").Value = "coyni.admin.tests.ProfilesTest,`ntestveifyApiBusinessDetails"
$ws.Range("I66").WrapText = $true
$ws.Range("I66").VerticalAlignment = -4160

# Apply the same cell styling pattern used throughout the sheet:
# wrap text for column H (Keywords1) on the new rows.
$ws.Range("H64:H66").WrapText = $true

$ws.Range("I66").Select()
